{"js": "// Replace each two-digit multiplication expression's text with its new value.\n// Each old expression is unique within the document, so an exact search +\n// full-text replace on the matched range is safe and order independent.\nconst replacements = [\n  [\"39\u00d757=\", \"62\u00d792=\"],\n  [\"75\u00d773=\", \"14\u00d714=\"],\n  [\"45\u00d715=\", \"15\u00d764=\"],\n  [\"38\u00d776=\", \"87\u00d788=\"],\n  [\"20\u00d742=\", \"90\u00d737=\"],\n  [\"96\u00d758=\", \"39\u00d759=\"],\n  [\"52\u00d780=\", \"98\u00d711=\"],\n  [\"92\u00d797=\", \"18\u00d727=\"],\n  [\"73\u00d715=\", \"38\u00d788=\"],\n  [\"77\u00d798=\", \"50\u00d738=\"],\n  [\"48\u00d746=\", \"71\u00d732=\"],\n  [\"69\u00d754=\", \"36\u00d717=\"],\n  [\"52\u00d720=\", \"32\u00d772=\"],\n  [\"13\u00d714=\", \"58\u00d763=\"],\n  [\"29\u00d728=\", \"52\u00d758=\"],\n  [\"92\u00d795=\", \"83\u00d716=\"],\n  [\"47\u00d787=\", \"60\u00d754=\"],\n  [\"89\u00d758=\", \"24\u00d725=\"],\n  [\"54\u00d789=\", \"70\u00d746=\"],\n  [\"86\u00d768=\", \"21\u00d774=\"],\n  [\"96\u00d762=\", \"76\u00d733=\"],\n  [\"37\u00d733=\", \"89\u00d717=\"],\n  [\"15\u00d727=\", \"87\u00d738=\"],\n  [\"17\u00d742=\", \"91\u00d716=\"],\n  [\"50\u00d719=\", \"25\u00d785=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression's text with its new value.\n# Each old expression is unique within the document, so Find/Replace on the\n# whole document body (wdReplaceAll) is safe and order independent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"39\u00d757=\", \"62\u00d792=\"),\n    @(\"75\u00d773=\", \"14\u00d714=\"),\n    @(\"45\u00d715=\", \"15\u00d764=\"),\n    @(\"38\u00d776=\", \"87\u00d788=\"),\n    @(\"20\u00d742=\", \"90\u00d737=\"),\n    @(\"96\u00d758=\", \"39\u00d759=\"),\n    @(\"52\u00d780=\", \"98\u00d711=\"),\n    @(\"92\u00d797=\", \"18\u00d727=\"),\n    @(\"73\u00d715=\", \"38\u00d788=\"),\n    @(\"77\u00d798=\", \"50\u00d738=\"),\n    @(\"48\u00d746=\", \"71\u00d732=\"),\n    @(\"69\u00d754=\", \"36\u00d717=\"),\n    @(\"52\u00d720=\", \"32\u00d772=\"),\n    @(\"13\u00d714=\", \"58\u00d763=\"),\n    @(\"29\u00d728=\", \"52\u00d758=\"),\n    @(\"92\u00d795=\", \"83\u00d716=\"),\n    @(\"47\u00d787=\", \"60\u00d754=\"),\n    @(\"89\u00d758=\", \"24\u00d725=\"),\n    @(\"54\u00d789=\", \"70\u00d746=\"),\n    @(\"86\u00d768=\", \"21\u00d774=\"),\n    @(\"96\u00d762=\", \"76\u00d733=\"),\n    @(\"37\u00d733=\", \"89\u00d717=\"),\n    @(\"15\u00d727=\", \"87\u00d738=\"),\n    @(\"17\u00d742=\", \"91\u00d716=\"),\n    @(\"50\u00d719=\", \"25\u00d785=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
